# "fixed the reservation system"
# Adds a new order record (row 17) to the "Order Data" sheet, representing
# a delivery order (Customer ID 4) with items "[1, 1, 2]" that is currently
# InProgress and not yet completed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New order row
$ws.Range("A17").Value = 18
$ws.Range("B17").Value = "delivery"
$ws.Range("C17").Value = "[1, 1, 2]"
$ws.Range("D17").Value = $false
$ws.Range("E17").Value = "InProgress"
$ws.Range("F17").Value = 4
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0

# Re-apply the (near) best-fit column widths for the data columns, matching
# the tiny width recalculation Excel performs whenever the backing data
# driving "best fit" columns changes.
$ws.Columns.Item(1).ColumnWidth = 7.833333333333334
$ws.Columns.Item(2).ColumnWidth = 10.333333333333332
$ws.Columns.Item(3).ColumnWidth = 30.833333333333336
$ws.Columns.Item(4).ColumnWidth = 21.833333333333336
$ws.Columns.Item(5).ColumnWidth = 13.333333333333332
$ws.Columns.Item(6).ColumnWidth = 11.333333333333332
